# Weekly update: insert two new daily price records for Membrillo
# (Mercado Mayorista Lo Valledor de Santiago) at the top of the data
# block, pushing the existing rows 90-134 down to 92-136.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right above the current row 90.
$ws.Rows.Item(90).Insert()
$ws.Rows.Item(90).Insert()

# New row 90
$ws.Cells.Item(90, 1).Value  = 6
$ws.Cells.Item(90, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(90, 3).Value  = "Metropolitana"
$ws.Cells.Item(90, 4).Value  = 45001
$ws.Cells.Item(90, 5).Value  = 13
$ws.Cells.Item(90, 6).Value  = "Fruta"
$ws.Cells.Item(90, 7).Value  = 100104
$ws.Cells.Item(90, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(90, 9).Value  = 100104003
$ws.Cells.Item(90, 10).Value = "Membrillo"
$ws.Cells.Item(90, 11).Value = "Champion"
$ws.Cells.Item(90, 12).Value = "Especial"
$ws.Cells.Item(90, 13).Value = 15
$ws.Cells.Item(90, 14).Value = 280000
$ws.Cells.Item(90, 15).Value = 280000
$ws.Cells.Item(90, 16).Value = 280000
$ws.Cells.Item(90, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(90, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(90, 19).Value = 622
$ws.Cells.Item(90, 20).Value = 450

# New row 91
$ws.Cells.Item(91, 1).Value  = 6
$ws.Cells.Item(91, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(91, 3).Value  = "Metropolitana"
$ws.Cells.Item(91, 4).Value  = 45001
$ws.Cells.Item(91, 5).Value  = 13
$ws.Cells.Item(91, 6).Value  = "Fruta"
$ws.Cells.Item(91, 7).Value  = 100104
$ws.Cells.Item(91, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(91, 9).Value  = 100104003
$ws.Cells.Item(91, 10).Value = "Membrillo"
$ws.Cells.Item(91, 11).Value = "Champion"
$ws.Cells.Item(91, 12).Value = "Primera"
$ws.Cells.Item(91, 13).Value = 20
$ws.Cells.Item(91, 14).Value = 250000
$ws.Cells.Item(91, 15).Value = 250000
$ws.Cells.Item(91, 16).Value = 250000
$ws.Cells.Item(91, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(91, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(91, 19).Value = 556
$ws.Cells.Item(91, 20).Value = 450

Write-Output "Inserted 2 rows; new dimension should be A1:T136"
